$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.67309999999999
$ws.Range("A6").Value = -22.69340000000001
$ws.Range("A7").Value = -22.00539999999999
$ws.Range("B7").Value = 4.935000000000002
$ws.Range("A8").Value = -22.35070000000001
$ws.Range("B11").Value = 5.538200000000002
$ws.Range("B12").Value = 4.772699999999998
$ws.Range("B15").Value = 5.119399999999995
$ws.Range("A16").Value = -21.51989999999997
$ws.Range("A20").Value = -23.04560000000001
$ws.Range("B20").Value = 4.811299999999996
$ws.Range("A21").Value = -22.29569999999999
$ws.Range("B21").Value = 5.591199999999998
$ws.Range("B22").Value = 10.1215
$ws.Range("B23").Value = 9.227500000000004
$ws.Range("A28").Value = -22.21239999999999
$ws.Range("A29").Value = -21.7216
$ws.Range("B29").Value = 5.081100000000001
$ws.Range("A30").Value = -21.82180000000001
$ws.Range("A32").Value = -21.26739999999998
$ws.Range("B34").Value = 8.985400000000004
$ws.Range("A40").Value = -19.74259999999999
$ws.Range("B42").Value = 9.759699999999999
$ws.Range("B43").Value = 5.468800000000001
$ws.Range("B44").Value = 5.4715
$ws.Range("B45").Value = 4.9286
$ws.Range("A46").Value = -22.11699999999999
$ws.Range("B46").Value = 5.321500000000002
$ws.Range("B50").Value = 4.494599999999997
$ws.Range("A51").Value = -22.20659999999999
$ws.Range("B51").Value = 5.480199999999996
$ws.Range("A52").Value = -22.0698
$ws.Range("A57").Value = -22.66560000000001
$ws.Range("B57").Value = 4.916899999999995
$ws.Range("A59").Value = -22.25590000000001
$ws.Range("A62").Value = -22.1767
$ws.Range("B65").Value = 5.338000000000001
$ws.Range("A66").Value = -21.486
$ws.Range("B66").Value = 4.965399999999995
$ws.Range("B67").Value = 5.125599999999999
$ws.Range("A73").Value = -20.27599999999999
$ws.Range("A74").Value = -21.92419999999999
$ws.Range("A77").Value = -20.23039999999998
$ws.Range("B79").Value = 9.612400000000003
$ws.Range("B84").Value = 5.347700000000001
$ws.Range("B87").Value = 4.905899999999999
$ws.Range("A92").Value = -21.49530000000002
$ws.Range("B92").Value = 4.696399999999995
$ws.Range("B97").Value = 6.317299999999998
$ws.Range("A100").Value = -22.1462
